# "Fixed one of the cape tables" - fill in the missing HR / CI values on the
# "cape Cox treat only" sheet (the table that had several blank/placeholder
# cells), correct a leftover "()" placeholder, widen column E so the longer
# CI strings are readable, and move the active-cell selection.

# E10 used to hold a stray empty "()" placeholder - fill in the real CI.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cape Cox treat only")

$ws.Range("E10").Value = "(0.966,1.383)"

# --- Column G (lo/hi 95 CI, styled like the rest of the table) -----------
$ws.Range("G8").Value = "(0.400,0.585)"
$ws.Range("G9").Value = "(0.354,0.481"
$ws.Range("G10").Value = "(0.981,1.402)"

# --- Column I (second lo/hi 95 CI) ----------------------------------------
$ws.Range("I8").Value = "(0.543,0.906)"
$ws.Range("I9").Value = "(0.470,0.748)"
$ws.Range("I10").Value = "X FILL THIS IN"

# --- Cosmetics -------------------------------------------------------------
# Widen column E (it now holds longer confidence-interval strings).
# ColumnWidth is in characters; the engine rounds to whole pixels, so feed it
# the (pre-rounding) character width that corresponds to the target stored
# width of 13.42578125.
$ws.Columns.Item(5).ColumnWidth = 12.592447916666666

# Move the kept-alive selection.
$ws.Range("K13").Select()
